$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained one more row (a new "2007" forecast year was inserted at the
# top of the data, shifting the existing 2008-2024 rows down by one) and every
# y_0_forecast / y_1_forecast value was recomputed. Rather than using a native
# row-insert (which in this runtime ends up synthesizing extra/blended cell
# styles), write the full, final 19-row table directly: this naturally shifts
# rows 2-18 down to 3-19 and adds the new row 19, reusing the existing styles.

# Give the brand-new last row (19) the same date-column formatting used by
# every other row in column A (same style as the previous last row, 18).
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row, A(date_of_forecast), B(y_0), C(y_0_forecast), D(y_1), E(y_1_forecast)
$data = @(
    @(2, 39400, 2007, 1.75539628881467, 2008, 0.2337905658324813),
    @(3, 39765, 2008, 2.213911448916162, 2009, 3.386383090739953),
    @(4, 40130, 2009, 2.533533936850563, 2010, 0.984293482975751),
    @(5, 40494, 2010, 2.088987486264915, 2011, 3.612753212925401),
    @(6, 40862, 2011, 1.212544822741002, 2012, 2.158838189283174),
    @(7, 41228, 2012, 1.196776590518644, 2013, 1.194058515117336),
    @(8, 41592, 2013, 0.4712609263772594, 2014, 1.409662779709819),
    @(9, 41957, 2014, 0.8783377572271434, 2015, 2.372074663906587),
    @(10, 42321, 2015, 2.29066283401107, 2016, 4.595879021798321),
    @(11, 42689, 2016, 4.109890522944348, 2017, 4.034919509273061),
    @(12, 43053, 2017, 1.336316831462692, 2018, 0.02883756256675252),
    @(13, 43418, 2018, 1.197912858979611, 2019, 0.9262553939922924),
    @(14, 43783, 2019, 1.727537197898665, 2020, 2.928189816005666),
    @(15, 44159, 2020, 3.647228437274408, 2021, 3.673004547855219),
    @(16, 44525, 2021, 2.777797690741424, 2022, 1.579011422502852),
    @(17, 44890, 2022, 0.6994919452575576, 2023, -2.087978868409623),
    @(18, 45254, 2023, -1.432689847121871, 2024, 0.1172571542027212),
    @(19, 45618, 2024, 2.033479419175133, 2025, 1.317145539573517)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
